$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: find the 1-based index of the first paragraph whose text matches
# a given substring.
# ---------------------------------------------------------------------------
function Find-ParagraphIndex($doc, [string]$needle) {
    $count = $doc.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $t = $doc.Paragraphs.Item($i).Range.Text
        if ($t -like "*$needle*") {
            return $i
        }
    }
    return -1
}

# ---------------------------------------------------------------------------
# 0) Drop the (currently unique) "_GoBack" bookmark that sits at the very
#    end of the document *before* we insert any new text. Word only ever
#    keeps one "_GoBack" bookmark alive -- the edit below re-creates it at
#    the new insertion point, mirroring where the author last typed.
# ---------------------------------------------------------------------------
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# ---------------------------------------------------------------------------
# 1) Insert a brand-new bullet ("Empty divs ... as placeholders for the
#    sound file") right before the "Start button to begin the game" bullet.
# ---------------------------------------------------------------------------
$startIdx = Find-ParagraphIndex $d "Start button to begin the game"
$hostPara = $d.Paragraphs.Item($startIdx)
$hostPara.Range.InsertParagraphBefore() | Out-Null

$newPara = $d.Paragraphs.Item($startIdx)

$newParaXml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="ListParagraph"/>
              <w:numPr>
                <w:ilvl w:val="2"/>
                <w:numId w:val="7"/>
              </w:numPr>
              <w:shd w:val="clear" w:color="auto" w:fill="FFFFFF"/>
              <w:spacing w:before="100" w:beforeAutospacing="1" w:after="100" w:afterAutospacing="1"/>
              <w:contextualSpacing w:val="0"/>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
                <w:iCs w:val="0"/>
                <w:color w:val="333333"/>
                <w:spacing w:val="3"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
            </w:pPr>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
                <w:iCs w:val="0"/>
                <w:color w:val="333333"/>
                <w:spacing w:val="3"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve">Empty </w:t>
            </w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
                <w:iCs w:val="0"/>
                <w:color w:val="333333"/>
                <w:spacing w:val="3"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>divs</w:t>
            </w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
                <w:iCs w:val="0"/>
                <w:color w:val="333333"/>
                <w:spacing w:val="3"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t xml:space="preserve"> </w:t>
            </w:r>
            <w:bookmarkStart w:id="0" w:name="_GoBack"/>
            <w:bookmarkEnd w:id="0"/>
            <w:r>
              <w:rPr>
                <w:rFonts w:ascii="Helvetica Neue" w:eastAsia="Times New Roman" w:hAnsi="Helvetica Neue" w:cs="Times New Roman"/>
                <w:iCs w:val="0"/>
                <w:color w:val="333333"/>
                <w:spacing w:val="3"/>
                <w:sz w:val="24"/>
                <w:szCs w:val="24"/>
              </w:rPr>
              <w:t>as placeholders for the sound file</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$newPara.Range.InsertXML($newParaXml) | Out-Null

# ---------------------------------------------------------------------------
# 2) "Click event button to start gam" + "e" -> one run reading
#    "Click event button to start game" (fixes the split-word typo).
# ---------------------------------------------------------------------------
$d.Content.Find.Execute("Click event button to start game", $false, $false, $false, $false, $false, $true, 1, $false, "Click event button to start game", 2) | Out-Null

# ---------------------------------------------------------------------------
# 3) Remove the stray empty paragraph (ind left=2160) that used to sit right
#    before the "Function will need to enable multiple color inputs..." bullet.
# ---------------------------------------------------------------------------
$funcIdx = Find-ParagraphIndex $d "Function will need to enable multiple color inputs"
$emptyPara = $d.Paragraphs.Item($funcIdx + 1)
if ($emptyPara.Range.Text.Trim().Length -eq 0) {
    $emptyPara.Range.Delete()
}

Write-Output "edit complete"
